# Add research paper presentation
# Duplicates the measurement table (columns Plik / Liczba klatek / czasy / laczny czas /
# sekunda filmu) further down the sheet in a condensed form (dropping the
# "Detekcja" / "Interpolacja" detail columns), plus two small helper formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- small helper formulas (rows 14 and 16) -------------------------------
$ws.Range("D14").Formula = "=E13/11"
$ws.Range("F16").Formula = "=C13/G13"

# --- condensed table header (row 23) ---------------------------------------
$ws.Range("A23").Value = "Plik"
$ws.Range("B23").Value = "Liczba klatek filmu"
$ws.Range("C23").Value = "Czas modelowania ruchu kamery"
$ws.Range("D23").Value = "Czas kalibracji"
$ws.Range("E23").Value = "Łącznie czas"
$ws.Range("F23").Value = "Czas przetwarzania sekundy materiału"

# --- condensed table data (rows 24-34), copied from the original table -----
$ws.Range("A24").Value = "baltyk_koszalin_02.mp4"
$ws.Range("A25").Value = "baltyk_koszalin_03_03.mp4"
$ws.Range("A26").Value = "baltyk_koszalin_04_04.mp4"
$ws.Range("A27").Value = "baltyk_koszalin_05_06.mp4"
$ws.Range("A28").Value = "baltyk_koszalin_06_07.mp4"
$ws.Range("A29").Value = "baltyk_koszalin_07_09.mp4"
$ws.Range("A30").Value = "baltyk_kotwica_1.mp4"
$ws.Range("A31").Value = "baltyk_starogard_1.mp4"
$ws.Range("A32").Value = "WDA_Kotwica_01.mp4"
$ws.Range("A33").Value = "ENG_POL_01_09.mp4"
$ws.Range("A34").Value = "BAR_SEV_01.mp4"

$ws.Range("B24").Value = 667
$ws.Range("B25").Value = 480
$ws.Range("B26").Value = 419
$ws.Range("B27").Value = 461
$ws.Range("B28").Value = 300
$ws.Range("B29").Value = 660
$ws.Range("B30").Value = 175
$ws.Range("B31").Value = 728
$ws.Range("B32").Value = 501
$ws.Range("B33").Value = 758
$ws.Range("B34").Value = 372

$ws.Range("C24").Value = 27.719
$ws.Range("C25").Value = 20.471
$ws.Range("C26").Value = 16.779
$ws.Range("C27").Value = 18.826
$ws.Range("C28").Value = 12.024
$ws.Range("C29").Value = 26.542
$ws.Range("C30").Value = 7.15
$ws.Range("C31").Value = 29.099
$ws.Range("C32").Value = 20.063
$ws.Range("C33").Value = 40.1403
$ws.Range("C34").Value = 22.4561

$ws.Range("D24").Value = 17.633
$ws.Range("D25").Value = 15.55
$ws.Range("D26").Value = 8.791
$ws.Range("D27").Value = 13.858
$ws.Range("D28").Value = 11.7
$ws.Range("D29").Value = 14.87
$ws.Range("D30").Value = 8.989
$ws.Range("D31").Value = 21.044
$ws.Range("D32").Value = 16.742
$ws.Range("D33").Value = 18.499
$ws.Range("D34").Value = 11.5269

# row 24 formulas are stand-alone (start of the block, not part of a shared group)
$ws.Range("E24").Formula = "=C24+D24"
$ws.Range("F24").Formula = "=E24/B24 * 24"

# rows 25:34 / 25:35 become shared formula groups, written as ranges so the
# relative references propagate row by row
$ws.Range("E25:E34").Formula = "=C25+D25"
$ws.Range("F25:F35").Formula = "=E25/B25 * 24"

# --- totals row (row 35) ----------------------------------------------------
$ws.Range("A35").Value = "Łącznie"
$ws.Range("B35").Formula = "=SUM(B24:B34)"
$ws.Range("C35:E35").Formula = "=SUM(C24:C34)"

# --- number formatting for the new table's numeric columns -----------------
$ws.Range("C24:F35").NumberFormat = "0.000"

# --- column widths for the new table ---------------------------------------
$ws.Range("C1:E1").EntireColumn.ColumnWidth = 9.57
$ws.Range("F1:F1").EntireColumn.ColumnWidth = 9.29

# --- restore selection ------------------------------------------------------
$ws.Range("F17").Select()
